# [Outlook] Clean up (#503)
# Updates the "Member ID (methods only)" column (C) from 2 -> 1 for a
# specific set of rows in the Snippets sheet, and resets the sheet's
# scroll/selection view back to the top (A2) instead of being left
# scrolled down near row 139/157.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snippets")

# Reset the view: scroll position back to the top, frozen header pane,
# and active cell/selection back to A2.
$ws.Activate()
$ws.Range("A2").Select()

# Rows whose column C value changes from 2 to 1.
$rows = @(29, 38, 54, 59, 64, 65, 70, 75, 80, 89, 90, 108, 113, 117, 119, 121, 123, 187, 188, 189, 190, 191, 210, 223, 229, 232, 234, 236, 240, 244, 247, 250, 251)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = 1
}
